$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 7374.5
$ws.Range("I4").Value = 8413.857
$ws.Range("J4").Value = 99
$ws.Range("K4").Value = 8413.857
$ws.Range("L4").Value = 99
$ws.Range("M4").Value = -8299.857
$ws.Range("N4").Value = -327
$ws.Range("H39").Value = 495.3125
$ws.Range("I39").Value = 200.07692
$ws.Range("K39").Value = 600.23076
$ws.Range("M39").Value = -304.23076
$ws.Range("H44").Value = 101000
$ws.Range("J44").Value = 101000
$ws.Range("L44").Value = 101000
$ws.Range("N44").Value = -101924
$ws.Range("H137").Value = 4445.769
$ws.Range("I137").Value = 2056
$ws.Range("J137").Value = 8269.4
$ws.Range("K137").Value = 6168
$ws.Range("L137").Value = 24808.2
$ws.Range("M137").Value = -3618
$ws.Range("N137").Value = -29908.2
$ws.Range("H138").Value = 2041.7307
$ws.Range("J138").Value = 2699.5151
$ws.Range("L138").Value = 8098.5453
$ws.Range("N138").Value = -18378.5453

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 338.44446
$ws.Range("I5").Value = 242.58333
$ws.Range("J5").Value = 530.1667
$ws.Range("K5").Value = 242.58333
$ws.Range("L5").Value = 530.1667
$ws.Range("M5").Value = -130.58333
$ws.Range("N5").Value = -754.1667
$ws.Range("H32").Value = 21745268
$ws.Range("I32").Value = 23814270
$ws.Range("J32").Value = 20749
$ws.Range("K32").Value = 23814270
$ws.Range("L32").Value = 20749
$ws.Range("M32").Value = -23813983
$ws.Range("N32").Value = -21323
$ws.Range("H70").Value = 100000
$ws.Range("J70").Value = 100000
$ws.Range("L70").Value = 100000
$ws.Range("N70").Value = -100540
$ws.Range("H73").Value = 100000
$ws.Range("J73").Value = 100000
$ws.Range("L73").Value = 100000
$ws.Range("N73").Value = -101872
$ws.Range("H74").Value = 20849124
$ws.Range("I74").Value = 41667500
$ws.Range("K74").Value = 41667500
$ws.Range("M74").Value = -41666626
$ws.Range("H77").Value = 20849124
$ws.Range("I77").Value = 41667500
$ws.Range("K77").Value = 208337500
$ws.Range("M77").Value = -208333132
$ws.Range("H102").Value = 18153.312
$ws.Range("I102").Value = 19263.533
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 19263.533
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = -17641.533
$ws.Range("N102").Value = -4744

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 338.44446
$ws.Range("I4").Value = 242.58333
$ws.Range("J4").Value = 530.1667
$ws.Range("K4").Value = 242.58333
$ws.Range("L4").Value = 530.1667
$ws.Range("M4").Value = -127.58333
$ws.Range("N4").Value = -760.1667
$ws.Range("H99").Value = 2434.1667
$ws.Range("I99").Value = 1941
$ws.Range("J99").Value = 4900
$ws.Range("K99").Value = 1941
$ws.Range("L99").Value = 4900
$ws.Range("M99").Value = -443
$ws.Range("N99").Value = -7896

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3115.9092
$ws.Range("J7").Value = 5637.3335
$ws.Range("L7").Value = 5637.3335
$ws.Range("N7").Value = -5863.3335
$ws.Range("H31").Value = 492459.12
$ws.Range("I31").Value = 8740.344999999999
$ws.Range("J31").Value = 1230766.8
$ws.Range("K31").Value = 8740.344999999999
$ws.Range("L31").Value = 1230766.8
$ws.Range("M31").Value = -8445.344999999999
$ws.Range("N31").Value = -1231356.8
$ws.Range("H34").Value = 492459.12
$ws.Range("I34").Value = 8740.344999999999
$ws.Range("J34").Value = 1230766.8
$ws.Range("K34").Value = 8740.344999999999
$ws.Range("L34").Value = 1230766.8
$ws.Range("M34").Value = -8538.344999999999
$ws.Range("N34").Value = -1231170.8
$ws.Range("H122").Value = 1346
$ws.Range("I122").Value = 1397
$ws.Range("J122").Value = 1244
$ws.Range("K122").Value = 4191
$ws.Range("L122").Value = 3732
$ws.Range("M122").Value = -1741
$ws.Range("N122").Value = -8632
$ws.Range("H132").Value = 1732.0416
$ws.Range("I132").Value = 1589.9565
$ws.Range("K132").Value = 4769.8695
$ws.Range("M132").Value = -2239.8695

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 190698860
$ws.Range("I75").Value = 1000000000
$ws.Range("J75").Value = 55815344
$ws.Range("K75").Value = 3000000000
$ws.Range("L75").Value = 167446032
$ws.Range("M75").Value = -2999999002
$ws.Range("N75").Value = -167448028
$ws.Range("H78").Value = 190698860
$ws.Range("I78").Value = 1000000000
$ws.Range("J78").Value = 55815344
$ws.Range("K78").Value = 9000000000
$ws.Range("L78").Value = 502338096
$ws.Range("M78").Value = -8999995008
$ws.Range("N78").Value = -502348080
$ws.Range("H134").Value = 4137.3477
$ws.Range("J134").Value = 8000
$ws.Range("L134").Value = 24000
$ws.Range("N134").Value = -34140
$ws.Range("H139").Value = 2631.5789
$ws.Range("I139").Value = 2500
$ws.Range("J139").Value = 2638.889
$ws.Range("K139").Value = 7500
$ws.Range("L139").Value = 7916.667
$ws.Range("M139").Value = -2360
$ws.Range("N139").Value = -18196.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 110872
$ws.Range("J69").Value = 110872
$ws.Range("L69").Value = 110872
$ws.Range("N69").Value = -112370
$ws.Range("H72").Value = 110872
$ws.Range("J72").Value = 110872
$ws.Range("L72").Value = 332616
$ws.Range("N72").Value = -340104
$ws.Range("H80").Value = 1749.75
$ws.Range("J80").Value = 999.5
$ws.Range("L80").Value = 999.5
$ws.Range("N80").Value = -2995.5
$ws.Range("H83").Value = 1749.75
$ws.Range("J83").Value = 999.5
$ws.Range("L83").Value = 4997.5
$ws.Range("N83").Value = -14981.5
$ws.Range("H102").Value = 5155.8335
$ws.Range("I102").Value = 3530.111
$ws.Range("J102").Value = 10033
$ws.Range("K102").Value = 3530.111
$ws.Range("L102").Value = 10033
$ws.Range("M102").Value = -1908.111
$ws.Range("N102").Value = -13277
$ws.Range("H126").Value = 3119.4375
$ws.Range("I126").Value = 2943.6
$ws.Range("J126").Value = 3412.5
$ws.Range("K126").Value = 8830.799999999999
$ws.Range("L126").Value = 10237.5
$ws.Range("M126").Value = -6360.799999999999
$ws.Range("N126").Value = -15177.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11170410
$ws.Range("I7").Value = 18183868
$ws.Range("K7").Value = 18183868
$ws.Range("M7").Value = -18183756
$ws.Range("H40").Value = 3064.658
$ws.Range("I40").Value = 2515.7932
$ws.Range("J40").Value = 4833.222
$ws.Range("K40").Value = 2515.7932
$ws.Range("L40").Value = 4833.222
$ws.Range("M40").Value = -2379.7932
$ws.Range("N40").Value = -5105.222
$ws.Range("H82").Value = 529.2
$ws.Range("I82").Value = 511.5
$ws.Range("J82").Value = 600
$ws.Range("K82").Value = 511.5
$ws.Range("L82").Value = 600
$ws.Range("M82").Value = -150.5
$ws.Range("N82").Value = -1322
$ws.Range("H85").Value = 529.2
$ws.Range("I85").Value = 511.5
$ws.Range("J85").Value = 600
$ws.Range("K85").Value = 511.5
$ws.Range("L85").Value = 600
$ws.Range("M85").Value = 736.5
$ws.Range("N85").Value = -3096
$ws.Range("H126").Value = 11170410
$ws.Range("I126").Value = 18183868
$ws.Range("K126").Value = 54551604
$ws.Range("M126").Value = -54549134
$ws.Range("H132").Value = 71213.55499999999
$ws.Range("I132").Value = 9113.5
$ws.Range("J132").Value = 288563.75
$ws.Range("K132").Value = 27340.5
$ws.Range("L132").Value = 865691.25
$ws.Range("M132").Value = -24810.5
$ws.Range("N132").Value = -870751.25
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9263.294
$ws.Range("I122").Value = 5271.4546
$ws.Range("J122").Value = 16581.666
$ws.Range("K122").Value = 15814.3638
$ws.Range("L122").Value = 49744.99800000001
$ws.Range("M122").Value = -13364.3638
$ws.Range("N122").Value = -54644.99800000001
$ws.Range("H126").Value = 5390.1665
$ws.Range("I126").Value = 4147.75
$ws.Range("K126").Value = 12443.25
$ws.Range("M126").Value = -9973.25
$ws.Range("H132").Value = 12421.611
$ws.Range("I132").Value = 1684.9286
$ws.Range("J132").Value = 50000
$ws.Range("K132").Value = 5054.7858
$ws.Range("L132").Value = 150000
$ws.Range("M132").Value = -2524.7858
$ws.Range("N132").Value = -155060
